# Weekly data update: insert a new weekly record at row 28 (Fecha 2021-12-17),
# pushing the existing rows 28..141 down to 29..142.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 28; rows 28..141 shift to 29..142.
$ws.Rows("28:28").Insert()

# Populate the new row 28 with the new weekly record.
$ws.Cells.Item(28, 1).Value = 8
$ws.Cells.Item(28, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28, 3).Value = "Coquimbo"
$ws.Cells.Item(28, 4).Value = 44547
$ws.Cells.Item(28, 5).Value = 4
$ws.Cells.Item(28, 6).Value = 100112037
$ws.Cells.Item(28, 7).Value = "Cebollín"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 3200
$ws.Cells.Item(28, 11).Value = 900
$ws.Cells.Item(28, 12).Value = 1000
$ws.Cells.Item(28, 13).Value = 950
$ws.Cells.Item(28, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(28, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(28, 16).Value = 158
$ws.Cells.Item(28, 17).Value = 6
$ws.Cells.Item(28, 18).Value = "Hortaliza"
